# Regenerate save_data: use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 3
    6 = 0
    7 = 2
    8 = 1
    9 = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 3
    18 = 0
    19 = 0
    20 = 2
    21 = 3
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 3
    27 = 0
    28 = 1
    29 = 2
    31 = 2
    32 = 1
    33 = 2
    34 = 2
    35 = 2
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 2
    41 = 3
    42 = 2
    43 = 2
    44 = 1
    45 = 2
    46 = 1
    47 = 0
    48 = 1
    49 = 1
    50 = 1
    51 = 2
    52 = 3
    53 = 1
    54 = 2
    55 = 4
    56 = 1
    57 = 3
    58 = 2
    59 = 1
    60 = 1
    61 = 2
    62 = 4
    63 = 1
    64 = 2
    65 = 3
    66 = 1
    68 = 2
    69 = 0
    70 = 0
    71 = 0
    72 = 2
    73 = 2
    74 = 3
    75 = 1
    76 = 1
    77 = 1
    78 = 2
    79 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Keys.Count) K values in column G"
